$wb = $excel.ActiveWorkbook

# --- Content edit -------------------------------------------------------
# Sheet "1_" (the 2nd tab) held a population-growth question about an
# initial population of "6 million". Update it to ask about "6 billion"
# (and clarify the answer should be expressed "in billions").
$ws1 = $wb.Worksheets.Item("1_")
$newQuestion = "If our initial population is 6 billion, and the population grows at a rate of" + [char]160 + "α=.03, what is the population after 20 years (in billions)?"
$ws1.Range("A1").Value = $newQuestion

# --- UI / selection state -------------------------------------------------
# Reflect the saved selection on sheet "1_" (cell B6 was selected there).
$ws1.Activate()
$ws1.Range("B6").Select()

# The workbook was last saved with sheet "3_" active and cell B3 selected.
$ws2 = $wb.Worksheets.Item("3_")
$ws2.Activate()
$ws2.Range("B3").Select()
